# users.xlsx template update:
#  - add a new "EndTime" column (inserted before the trailing "don't remove"
#    note column, which shifts right from L to M)
#  - give the sample row a password example and an EndTime example
#  - fix the sample userId value (ali -> alit)
#  - move the active-cell selection
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert a new column at L; this pushes the existing "don't remove" note
# column (and the helper text under it) from L to M, preserving their
# values/styles automatically.
$null = $ws.Columns("L").Insert()

# New header for the inserted "EndTime" column.
$ws.Range("L1").Value = "EndTime"

# The "password" column already existed (K); it previously had no sample
# value under the header - give it one.
$ws.Range("K2").Value = "#Parrsoo2020#"

# Sample value for the new EndTime column.
$ws.Range("L2").Value = "1400/10/20 13:13:13.259"

# Fix the sample userId value.
$ws.Range("A2").Value = "alit"

# Best-effort column widths for the new/shifted columns (K, L, M).
$ws.Columns("K").ColumnWidth = 11.71
$ws.Columns("L").ColumnWidth = 20.86
$ws.Columns("M").ColumnWidth = 44.86

# Move the selected cell.
$null = $ws.Range("D10").Select()
